$wb = $excel.ActiveWorkbook

# Add the new worksheet "Sheet1" after the last existing sheet (ProductSlider)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet1"

# Header row (row 1) - DataSet / UserName / Password / ... / search
$ws.Range("A1").Value = "DataSet"
$ws.Range("B1").Value = "UserName"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "FirstName"
$ws.Range("E1").Value = "LastName"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "SegmentName"
$ws.Range("H1").Value = "taxclass"
$ws.Range("I1").Value = "website"
$ws.Range("J1").Value = "ApplyTo"
$ws.Range("K1").Value = "AssignedStatus"
$ws.Range("L1").Value = "Description"
$ws.Range("M1").Value = "customergroupname"
$ws.Range("N1").Value = "Updatecustomergroupname"
$ws.Range("O1").Value = "Adminsuccessmesaage"
$ws.Range("P1").Value = "errormessage"
$ws.Range("Q1").Value = "Tiletext"
$ws.Range("R1").Value = "textbutton"
$ws.Range("S1").Value = "component"
$ws.Range("T1").Value = "attribute"
$ws.Range("U1").Value = "CategorySelect"
$ws.Range("V1").Value = "productnames"
$ws.Range("W1").Value = "Filter"
$ws.Range("X1").Value = "Condition"
$ws.Range("Y1").Value = "VideoURL"
$ws.Range("Z1").Value = "CardTitle"
$ws.Range("AA1").Value = "Author"
$ws.Range("AB1").Value = "Rating"
$ws.Range("AC1").Value = "Quote"
$ws.Range("AD1").Value = "image"
$ws.Range("AE1").Value = "Buttontext"
$ws.Range("AF1").Value = "Buttontype"
$ws.Range("AG1").Value = "ButtonLinknavigation"
$ws.Range("AH1").Value = "ButtonlinkURL"
$ws.Range("AI1").Value = "Buttonlinkcategory"
$ws.Range("AJ1").Value = "Buttonlinkproduct"
$ws.Range("AK1").Value = "Buttonlinkpage"
$ws.Range("AL1").Value = "Categorydisplay"
$ws.Range("AM1").Value = "No.ofproductsdisplay"
$ws.Range("AN1").Value = "productcategory"
$ws.Range("AO1").Value = "Editpagetitle"
$ws.Range("AP1").Value = "datacontenttype"
$ws.Range("AQ1").Value = "Deletcomponent"
$ws.Range("AR1").Value = "headingtype"
$ws.Range("AS1").Value = " Description Type"
$ws.Range("AT1").Value = "alterantivetext"
$ws.Range("AU1").Value = "titleaatribute"
$ws.Range("AV1").Value = "Customergrouppagetitle"
$ws.Range("AW1").Value = "title"
$ws.Range("AX1").Value = "subtitle"
$ws.Range("AY1").Value = "MYAccountlinks"
$ws.Range("AZ1").Value = "Website"
$ws.Range("BA1").Value = "Group"
$ws.Range("BB1").Value = "DOB"
$ws.Range("BC1").Value = "State"
$ws.Range("BD1").Value = "Gender"
$ws.Range("BE1").Value = "Welcome Email"
$ws.Range("BF1").Value = "ProDeal"
$ws.Range("BG1").Value = "Acceptdate"
$ws.Range("BH1").Value = "IsEnabled"
$ws.Range("BI1").Value = "Language"
$ws.Range("BJ1").Value = "Street"
$ws.Range("BK1").Value = "City"
$ws.Range("BL1").Value = "Postcode"
$ws.Range("BM1").Value = "Phonenumber"
$ws.Range("BN1").Value = "Country"
$ws.Range("BO1").Value = "Region"
$ws.Range("BP1").Value = "ordernumber"
$ws.Range("BQ1").Value = "SKU"
$ws.Range("BR1").Value = "Productname"
$ws.Range("BS1").Value = "Style"
$ws.Range("BT1").Value = "Alignment"
$ws.Range("BU1").Value = "Color"
$ws.Range("BV1").Value = "Backgroud position"
$ws.Range("BW1").Value = "Attachment"
$ws.Range("BX1").Value = "pageTitle"
$ws.Range("BY1").Value = "URL"
$ws.Range("BZ1").Value = "preprodURL"
$ws.Range("CA1").Value = "mobilelayout"
$ws.Range("CB1").Value = "heading"
$ws.Range("CC1").Value = "TextColor"
$ws.Range("CD1").Value = "CTA Type"
$ws.Range("CE1").Value = "CTA Link"
$ws.Range("CF1").Value = "CTAText"
$ws.Range("CG1").Value = "CTAurl"
$ws.Range("CH1").Value = "categoryname"
$ws.Range("CI1").Value = "mrgtop"
$ws.Range("CJ1").Value = "mrgright"
$ws.Range("CK1").Value = "mrgbottom"
$ws.Range("CL1").Value = "mrgleft"
$ws.Range("CM1").Value = "paddingtop"
$ws.Range("CN1").Value = "paddingright"
$ws.Range("CO1").Value = "paddingbottom"
$ws.Range("CP1").Value = "paddingleft"
$ws.Range("CQ1").Value = "CSSclasses"
$ws.Range("CR1").Value = "Author"
$ws.Range("CS1").Value = "Quote"
$ws.Range("CT1").Value = "SubTitle"
$ws.Range("CU1").Value = "ChooseCondition"
$ws.Range("CV1").Value = "FF"
$ws.Range("CW1").Value = "Updateproductname"
$ws.Range("CX1").Value = "Price"
$ws.Range("CY1").Value = "Stock Status"
$ws.Range("CZ1").Value = "Updateprice"
$ws.Range("DA1").Value = "updateStock Status"
$ws.Range("DB1").Value = "search"

# Row 1 uses the existing yellow header style (same as the other sheets row 1)
$ws.Range("A1:DB1").Interior.Color = 65535

# Row 2 - AccountDetails / Productupdate data
$ws.Range("A2").Value = "AccountDetails"
$ws.Range("B2").Value = "mkoppanadam@helenoftroy.com"
$ws.Range("C2").Value = "Amtlmcflmipq1!"
$ws.Range("F2").Value = "mkoppanadam@helenoftroy.com"
$ws.Range("AW2").Value = "Dashboard / Magento Admin"
$ws.Range("BR2").Value = "QATEST product"
$ws.Range("BX2").Value = "Home Page "
$ws.Range("BY2").Value = "https://mcloud-na-stage.oxo.com/"
$ws.Range("BZ2").Value = "https://mcloud-na-preprod.oxo.com/"
$ws.Range("CW2").Value = "10 QATEST product "
$ws.Range("CY2").Value = "Out of Stock"
$ws.Range("DA2").Value = "In Stock"
$ws.Range("DB2").Value = "Lowest Price"
$ws.Range("CX2").Value = 5
$ws.Range("CZ2").Value = 30

# Row 3
$ws.Range("A3").Value = "Productupdate"
$ws.Range("AL3").Value = "Bottles,POP Containers"
$ws.Range("BR3").Value = "QATEST product"
$ws.Range("BY3").Value = "https://mcloud-na-stage.hydroflask.com//"
$ws.Range("BZ3").Value = "https://mcloud-na-preprod.hydroflask.com/"
$ws.Range("CW3").Value = "10 QATEST product "
$ws.Range("CY3").Value = "Out of Stock"
$ws.Range("DA3").Value = "In Stock"
$ws.Range("DB3").Value = "Lowest Price"
$ws.Range("CX3").Value = 5
$ws.Range("CZ3").Value = 30

# Row 4
$ws.Range("A4").Value = "Address"
$ws.Range("AW4").Value = "Home Page "
$ws.Range("BX4").Value = "Home Page (Hydroflask)"
$ws.Range("BY4").Value = "https://mcloud-na-stage.hydroflask.com//"
$ws.Range("BZ4").Value = "https://mcloud-na-preprod.hydroflask.com/"

# Selection matches target sheetView
$ws.Range("D5").Select()
